$wb = $excel.ActiveWorkbook

# Sheet "Plate8": update Location codes for rows 18-21
$ws2 = $wb.Worksheets.Item("Plate8")
$ws2.Range("A18").Value = "r07c09"
$ws2.Range("A19").Value = "r07c11"
$ws2.Range("A20").Value = "r08c10"
$ws2.Range("A21").Value = "r08c12"

# Sheet "HUVEC_Control": update Treatment for row 10 (Location r03c03)
$ws3 = $wb.Worksheets.Item("HUVEC_Control")
$ws3.Range("B10").Value = "untreated"
